$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-log rows for the "Design" and "Tilaus luokka" tasks, continuing the
# small G/H/I side-table (task name / start time / end time) started at row 7.
$ws.Range("G10").Value = "Design"
$ws.Range("H10").Value = 0.79166666666666663
$ws.Range("H10").NumberFormat = "h:mm"
$ws.Range("I10").Value = 0.80555555555555547
$ws.Range("I10").NumberFormat = "h:mm"

$ws.Range("G11").Value = "Tilaus luokka"
$ws.Range("H11").Value = 0.80555555555555547
$ws.Range("H11").NumberFormat = "h:mm"
$ws.Range("I11").Value = 0.84027777777777779
$ws.Range("I11").NumberFormat = "h:mm"

# Move the active cell/selection the same way the author's session ended up.
$ws.Range("I12").Select() | Out-Null
